# Adds a new "Bayanic Calendar Rule:" bullet item to the editing-rules
# list, right after the "Readability rule" bullet and before the
# trailing empty paragraph.

$d = $word.ActiveDocument

# Locate the "Readability rule" paragraph (the last populated bullet,
# ending in "...break up sentences. ") by searching for unique text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Do not use any hyphens or dashes to break up sentences.*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the Readability rule paragraph"
}

# Insert a new empty paragraph right after it; this new paragraph
# inherits the numbered-list formatting (ilvl 0 / numId 3) already
# used throughout the list.
$target.Range.InsertParagraphAfter()

# Re-find it: it is the paragraph immediately after $target now.
$newIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$r = $newPara.Range

# Build the paragraph's exact OOXML (including the numbering
# properties, bold "Bayanic Calendar Rule:" heading, and the
# proofErr spell-check markers Word leaves around the unrecognized
# word "Bayanic") and inject it via InsertXML so formatting and run
# boundaries match exactly.
$innerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Bayanic</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Calendar Rule:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> If any dates occur, make the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Bayanic</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> date the primary date and the Gregorian date in parenthesis. The </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Bayanic</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Calendar Year 0 starts May 23, 1844 with Year 1 starting on the spring equinox of 1845. </w:t></w:r>' +
    '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $innerXml +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
